# Regenerate merged AHB files
# - rename the "_old"/"_new" column-suffix headers to the actual form-version
#   tags ("_FV2210" / "_FV2304")
# - turn the sheet's used range into a real Excel Table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. rename header cells (A1:J1 = "old"/FV2210 side, L1:U1 = "new"/FV2304 side) ---
$oldHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $leftCol  = $i + 1        # A..J
    $rightCol = $i + 12       # L..U
    $ws.Cells.Item(1, $leftCol).Value  = $oldHeaders[$i] + "_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = $oldHeaders[$i] + "_FV2304"
}

# --- 2. convert the data range into a table ---
$rng = $ws.Range("A1:U63")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
